$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text changes from "Ready for handoff" to
#     "Handback transform failed" for the 06e0a089... (.md) row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: same status change, plus a new error-detail cell (L3)
#     describing the handback/handoff file name mismatch.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("L3").Value = "Handback file name: sidkmqan.dr3 is different with handoff file name: 06e0a089-2a38-45c5-882f-32fe4bcf57e4.9b6de8173dc4b16303c2e5288ecbed303d55f615.zh-cn."

# --- de-de sheet: same status change, plus a new error-detail cell (L3)
#     describing the handback/handoff file name mismatch.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("L3").Value = "Handback file name: sidkmqan.dr3 is different with handoff file name: 06e0a089-2a38-45c5-882f-32fe4bcf57e4.9b6de8173dc4b16303c2e5288ecbed303d55f615.de-de."
